# Update the date/title line and the multiplication problems in the
# three-digit x one-digit multiplication worksheet.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-12-08 Monday" "2025-12-09 Tuesday"

Replace-Text "263×5=" "889×8="
Replace-Text "461×2=" "747×4="
Replace-Text "767×5=" "776×3="
Replace-Text "751×5=" "720×8="
Replace-Text "794×4=" "302×3="

Replace-Text "858×6=" "454×6="
Replace-Text "487×9=" "965×5="
Replace-Text "562×6=" "184×6="
Replace-Text "804×8=" "819×6="
Replace-Text "415×4=" "895×9="

Replace-Text "304×7=" "754×3="
Replace-Text "836×5=" "589×8="
Replace-Text "460×4=" "200×6="
Replace-Text "857×9=" "321×2="
Replace-Text "651×9=" "361×5="

Replace-Text "240×4=" "280×6="
Replace-Text "649×9=" "490×3="
Replace-Text "390×4=" "677×4="
Replace-Text "112×4=" "749×3="
Replace-Text "460×2=" "424×2="

Replace-Text "538×6=" "306×6="
Replace-Text "831×9=" "143×2="
Replace-Text "176×3=" "594×4="
Replace-Text "476×9=" "175×9="
Replace-Text "751×3=" "305×4="
